# Applies the "finished poster" edit:
#  - re-wraps three paragraphs of Text Box 19 (Evaluation Metrics) into
#    multiple runs, tweaking a couple of words along the way
#  - re-wraps paragraphs of the Motivation text box (Text Box 6, id=23)
#    into multiple runs, tweaking a couple of words, and nudges/widens
#    its bounding box
#  - re-wraps the "Currently several technologies..." paragraph of the
#    Cost Reduction Techniques text box (Text Box 6, id=27) into
#    multiple runs (wording unchanged)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Split-ParagraphRuns {
    param($Paragraph, $Chunks)

    $sz = $Paragraph.Font.Size
    $full = [string]::Join("", $Chunks)

    # The engine diffs a freshly-assigned .Text against the paragraph's
    # previous contents and auto-splits off any matching leading/trailing
    # run, which would otherwise throw off our chunk boundaries below.
    # Stomp the paragraph with text that shares no characters with any
    # of our target wording first so that assigning the real text
    # collapses cleanly down to a single run.
    $Paragraph.Text = "~~~"
    $Paragraph.Text = $full

    # Now force a run boundary after every chunk (except the last) by
    # re-asserting the (unchanged) font size on the remaining tail of
    # the text -- PowerPoint splits off a fresh run for the re-touched
    # range without adding any extra formatting attributes.
    $pos = 1
    for ($i = 0; $i -lt $Chunks.Length - 1; $i++) {
        $pos = $pos + $Chunks[$i].Length
        $remaining = $full.Length - $pos + 1
        $Paragraph.Characters($pos, $remaining).Font.Size = $sz
    }
}

# --- Text Box 19 ("Evaluation Metrics and Results" results box) ---
$shp11 = $s.Shapes.Item(11)
$tr11 = $shp11.TextFrame.TextRange

Split-ParagraphRuns $tr11.Paragraphs(1, 1) @(
    "When evaluating the efficacy of the research, several groups of images were developed containing unique and duplicate images in various sizes. Metrics ",
    "recorded include processing ",
    "time, detection rates, and storage requirements."
)

Split-ParagraphRuns $tr11.Paragraphs(3, 1) @(
    "Processing times increased, but ",
    "remained ",
    "less ",
    "than two seconds in observed scenarios, as seen in Figure 3."
)

# --- Text Box 6, id 23 (Motivation bullets) ---
$shp18 = $s.Shapes.Item(18)

# Nudge/widen the bounding box (EMU -> points, 12700 EMU per point)
$shp18.Left = 914400 / 12700
$shp18.Top = 18973800 / 12700
$shp18.Width = 13411200 / 12700
$shp18.Height = 4648200 / 12700

$tr18 = $shp18.TextFrame.TextRange

Split-ParagraphRuns $tr18.Paragraphs(1, 1) @(
    "As of May 2013, nearly 500 ",
    "million images were shared ",
    "each ",
    "day. This is ",
    "expected to double by May 2014. [1]"
)

Split-ParagraphRuns $tr18.Paragraphs(2, 1) @(
    "Approximately 20% of this ",
    "data ",
    "is estimated to be duplicate. [2]"
)

Split-ParagraphRuns $tr18.Paragraphs(3, 1) @(
    "By eliminating this ",
    "duplicate data",
    ", companies can ",
    "save roughly $",
    "1.8 million annually."
)

# --- Text Box 6, id 27 (Cost Reduction Techniques intro) ---
$shp22 = $s.Shapes.Item(22)
$tr22 = $shp22.TextFrame.TextRange

Split-ParagraphRuns $tr22.Paragraphs(1, 1) @(
    "Currently several technologies are used to reduce the costs associated with ",
    "storing the ",
    "shared data including:"
)
